$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unify SQL syntax casing (row 6 / PreviousMedBrand) ---
$ws.Range("C6").Value = "SELECT * FROM Brands"

# --- Row 7 (PreviousMedProduct) ---
$ws.Range("C7").Value = "SELECT * FROM Products WHERE brand_id == [prev_med_brand]"
$ws.Range("D7").Value = "Please choose your product from the med brand"

# New Action cell E7: rich text "GO(" + "SeeProduct)"
$ws.Range("E7").Value = "GO(SeeProduct)"
$e7b = $ws.Range("E7").Characters(4, 11)
$e7b.Font.Name = "Calibri"
$e7b.Font.Size = 11
$e7b.Font.Color = 0

$ws.Range("F7").Value = "[med]"

$ws.Rows.Item(7).RowHeight = 14.9

# --- Row 8 (SeeProduct) ---
$ws.Range("C8").Value = "SELECT * FROM Products WHERE id == [med]"

# Action cell E8: rich text "GO(" + "FreeText)"
$ws.Range("E8").Value = "GO(FreeText)"
$e8b = $ws.Range("E8").Characters(4, 9)
$e8b.Font.Name = "Calibri"
$e8b.Font.Size = 11
$e8b.Font.Color = 0

$ws.Rows.Item(8).RowHeight = 14.9

# --- Row 9 (FreeText) ---
$ws.Range("E9").Value = "FINISH()"

# --- Column E width shrinks ---
$ws.Columns.Item(5).ColumnWidth = 31.86

# --- Selection / view moves from E4 to D30 ---
$ws.Range("D30").Select() | Out-Null
